$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 264 ---
# Copy formats from the row above so the new cells reuse the existing
# date / 2-decimal number styles instead of creating new ones.
$ws.Cells.Item(263, 1).Copy()
$ws.Cells.Item(264, 1).PasteSpecial(-4122)
$ws.Cells.Item(264, 1).Value = 45307.541666666664

$ws.Cells.Item(263, 2).Copy()
$ws.Cells.Item(264, 2).PasteSpecial(-4122)
$ws.Cells.Item(264, 2).Value = 1.92

$ws.Cells.Item(264, 4).Value = "Ice on gauge, difficult to read"

# --- Row 265 ---
$ws.Cells.Item(264, 1).Copy()
$ws.Cells.Item(265, 1).PasteSpecial(-4122)
$ws.Cells.Item(265, 1).Value = 45308.479166666664

$ws.Cells.Item(264, 2).Copy()
$ws.Cells.Item(265, 2).PasteSpecial(-4122)
$ws.Cells.Item(265, 2).Value = 1.91

$ws.Cells.Item(265, 3).Value = 0.02

$ws.Cells.Item(265, 4).Value = "Snowfall SWE"

# --- Row 266 ---
$ws.Cells.Item(265, 1).Copy()
$ws.Cells.Item(266, 1).PasteSpecial(-4122)
$ws.Cells.Item(266, 1).Value = 45309.510416666664

$ws.Cells.Item(265, 2).Copy()
$ws.Cells.Item(266, 2).PasteSpecial(-4122)
$ws.Cells.Item(266, 2).Value = 1.9

# --- Row 267 ---
$ws.Cells.Item(266, 1).Copy()
$ws.Cells.Item(267, 1).PasteSpecial(-4122)
$ws.Cells.Item(267, 1).Value = 45310.510416666664

$ws.Cells.Item(266, 2).Copy()
$ws.Cells.Item(267, 2).PasteSpecial(-4122)
$ws.Cells.Item(267, 2).Value = 1.9

$ws.Cells.Item(267, 3).Value = 0.02

$ws.Cells.Item(267, 4).Value = "Snowfall SWE"

$excel.CutCopyMode = 0

# Update selection to match new bottom of sheet
$ws.Range("C267:D267").Select()
